$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "25.10", "333.39") stay as text
# by pre-formatting the target cells as Text before assigning values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.808.24"
$ws.Range("E2").Value = "  +2.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.114.85"
$ws.Range("E3").Value = "  +6.68%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.39"
$ws.Range("E5").Value = "  +3.24%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5317"
$ws.Range("E7").Value = "  +3.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4375"
$ws.Range("E8").Value = "  +6.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09005"
$ws.Range("E9").Value = "  +6.99%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.28"
$ws.Range("E10").Value = "  +8.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.179"
$ws.Range("E11").Value = "  +4.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.10"
$ws.Range("E12").Value = "  +4.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.109.58"
$ws.Range("E13").Value = "  +6.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.755"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.829"
$ws.Range("E15").Value = "  +5.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.49"
$ws.Range("E16").Value = "  +4.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001128"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("E21").Value = "  +0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.342"
$ws.Range("E22").Value = "  +4.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.865.21"
$ws.Range("E23").Value = "  +1.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.37"
$ws.Range("E24").Value = "  +7.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.355.11"
$ws.Range("E25").Value = "  +6.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.270"
$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.579"
$ws.Range("E28").Value = "  +8.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.96"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.35"
$ws.Range("E30").Value = "  +2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.173"
$ws.Range("E31").Value = "  +3.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1080"
$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.223"
$ws.Range("E33").Value = "  +3.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.016"
$ws.Range("E34").Value = "  +5.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.548"
$ws.Range("E35").Value = "  +17.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02617"
$ws.Range("E36").Value = "  +5.78%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.537"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.88"
$ws.Range("E38").Value = "  +9.39%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06737"
$ws.Range("E39").Value = "  +3.68%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.525"
$ws.Range("E40").Value = "  +6.50%  "

$ws.Range("E41").Value = "  +5.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6847"
$ws.Range("E42").Value = "  +4.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6439"
$ws.Range("E44").Value = "  +5.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("E46").Value = "  +3.92%  "

$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("E49").Value = "  +4.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.93"
$ws.Range("E50").Value = "  +4.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.78"
$ws.Range("E51").Value = "  -2.30%  "
